$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "dnasr281@gmail.com, System"
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "85.2%"
$ws.Range("S16").NumberFormat = "@"
$ws.Range("S16").Value = "79.0%"
$ws.Range("G17").Value = "dnasr281@gmail.com, System"
$ws.Range("S17").NumberFormat = "@"
$ws.Range("S17").Value = "68.4%"
$ws.Range("M18").Value = 21
$ws.Range("S18").NumberFormat = "@"
$ws.Range("S18").Value = "92.9%"
$ws.Range("S21").NumberFormat = "@"
$ws.Range("S21").Value = "97.4%"
$ws.Range("S22").NumberFormat = "@"
$ws.Range("S22").Value = "91.7%"
$ws.Range("S23").NumberFormat = "@"
$ws.Range("S23").Value = "88.3%"
$ws.Range("G24").Value = "dnasr281@gmail.com, System"
$ws.Range("H24").Value = "20/31"
$ws.Range("S24").NumberFormat = "@"
$ws.Range("S24").Value = "74.1%"
$ws.Range("H25").Value = "16/31"
$ws.Range("S25").NumberFormat = "@"
$ws.Range("S25").Value = "80.2%"
$ws.Range("S26").NumberFormat = "@"
$ws.Range("S26").Value = "73.3%"
$ws.Range("G39").Value = "dnasr281@gmail.com, System"
$ws.Range("G46").Value = "dnasr281@gmail.com, System"
$ws.Range("H46").Value = "7/19"
$ws.Range("H47").Value = "7/19"
$ws.Range("G61").Value = "dnasr281@gmail.com, System"
$ws.Range("G68").Value = "dnasr281@gmail.com, System"
$ws.Range("H68").Value = "18/21"
$ws.Range("H69").Value = "18/21"
$ws.Range("H70").Value = "21/21"
$ws.Range("H71").Value = "0/21"
$ws.Range("H72").Value = "0/21"
$ws.Range("H73").Value = "0/21"
$ws.Range("H74").Value = "0/21"
$ws.Range("H75").Value = "0/21"
$ws.Range("H76").Value = "0/21"
$ws.Range("H77").Value = "0/21"
$ws.Range("H78").Value = "0/21"
$ws.Range("H79").Value = "0/21"
$ws.Range("H80").Value = "0/21"
$ws.Range("H81").Value = "0/21"
$ws.Range("H82").Value = "0/21"
$ws.Range("G83").Value = "dnasr281@gmail.com, System"
$ws.Range("H83").Value = "21/21"
$ws.Range("H84").Value = "0/21"
$ws.Range("H85").Value = "0/21"
$ws.Range("H86").Value = "0/21"
$ws.Range("H87").Value = "0/21"
$ws.Range("H88").Value = "0/21"
$ws.Range("H89").Value = "0/21"
$ws.Range("G90").Value = "dnasr281@gmail.com, System"
$ws.Range("G105").Value = "dnasr281@gmail.com, System"
$ws.Range("G112").Value = "dnasr281@gmail.com, System"
$ws.Range("G127").Value = "dnasr281@gmail.com, System"
$ws.Range("G134").Value = "dnasr281@gmail.com, System"
$ws.Range("H135").Value = "26/29"
$ws.Range("G149").Value = "dnasr281@gmail.com, System"
$ws.Range("G156").Value = "dnasr281@gmail.com, System"
$ws.Range("H156").Value = "32/33"
$ws.Range("G171").Value = "dnasr281@gmail.com, System"
$ws.Range("G178").Value = "dnasr281@gmail.com, System"
$ws.Range("H178").Value = "24/30"
$ws.Range("G193").Value = "dnasr281@gmail.com, System"
$ws.Range("G200").Value = "dnasr281@gmail.com, System"
$ws.Range("H200").Value = "14/27"
$ws.Range("G215").Value = "dnasr281@gmail.com, System"
$ws.Range("G222").Value = "dnasr281@gmail.com, System"
$ws.Range("H222").Value = "20/29"
$ws.Range("G237").Value = "dnasr281@gmail.com, System"
$ws.Range("D244").NumberFormat = "@"
$ws.Range("D244").Value = "1"
$ws.Range("G244").Value = "dnasr281@gmail.com, System"
$ws.Range("H244").Value = "15/30"
$ws.Range("D245").NumberFormat = "@"
$ws.Range("D245").Value = "2"
$ws.Range("H245").Value = "13/30"
$ws.Range("D246").NumberFormat = "@"
$ws.Range("D246").Value = "3"
$ws.Range("D247").NumberFormat = "@"
$ws.Range("D247").Value = "4"
$ws.Range("D248").NumberFormat = "@"
$ws.Range("D248").Value = "5"
$ws.Range("D249").NumberFormat = "@"
$ws.Range("D249").Value = "6"
$ws.Range("D250").NumberFormat = "@"
$ws.Range("D250").Value = "7"
$ws.Range("D251").NumberFormat = "@"
$ws.Range("D251").Value = "8"
$ws.Range("D252").NumberFormat = "@"
$ws.Range("D252").Value = "9"
$ws.Range("D253").NumberFormat = "@"
$ws.Range("D253").Value = "10"
$ws.Range("D254").NumberFormat = "@"
$ws.Range("D254").Value = "11"
$ws.Range("D255").NumberFormat = "@"
$ws.Range("D255").Value = "12"
$ws.Range("D256").NumberFormat = "@"
$ws.Range("D256").Value = "13"
$ws.Range("D257").NumberFormat = "@"
$ws.Range("D257").Value = "14"
$ws.Range("D258").NumberFormat = "@"
$ws.Range("D258").Value = "15"
$ws.Range("D259").NumberFormat = "@"
$ws.Range("D259").Value = "16"
$ws.Range("G259").Value = "dnasr281@gmail.com, System"
$ws.Range("D260").NumberFormat = "@"
$ws.Range("D260").Value = "17"
$ws.Range("D261").NumberFormat = "@"
$ws.Range("D261").Value = "18"
$ws.Range("D262").NumberFormat = "@"
$ws.Range("D262").Value = "19"
$ws.Range("D263").NumberFormat = "@"
$ws.Range("D263").Value = "20"
$ws.Range("D264").NumberFormat = "@"
$ws.Range("D264").Value = "21"
$ws.Range("D265").NumberFormat = "@"
$ws.Range("D265").Value = "22"
